$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.510.91"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.984.33"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.31%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.37"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.34"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.353"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0721"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.81%  "
$ws.Range("E12").Value = "  -6.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.881"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.267.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.975.94"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.422.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.91"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0829"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -9.57%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -4.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.80%  "
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("E31").Value = "  -4.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.76"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -9.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0585"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0887"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +8.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -10.37%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.18"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -11.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.82"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("E40").Value = "  -7.93%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0207"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.377.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0874"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -8.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.33"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.73"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.82%  "
